$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 11.2
$ws.Range("I11").Value = 11.2
$ws.Range("K11").Value = 11.2
$ws.Range("M11").Value = 128.8
$ws.Range("H52").Value = 275.8
$ws.Range("I52").Value = 199.66667
$ws.Range("J52").Value = 390
$ws.Range("K52").Value = 599.00001
$ws.Range("L52").Value = 1170
$ws.Range("M52").Value = -439.00001
$ws.Range("N52").Value = -1490
$ws.Range("H112").Value = 1819.0834
$ws.Range("J112").Value = 1416.125
$ws.Range("L112").Value = 4248.375
$ws.Range("N112").Value = -6464.375
$ws.Range("H116").Value = 2129.1667
$ws.Range("I116").Value = 592.5
$ws.Range("J116").Value = 2897.5
$ws.Range("K116").Value = 592.5
$ws.Range("L116").Value = 2897.5
$ws.Range("M116").Value = 2849.5
$ws.Range("N116").Value = -9781.5
$ws.Range("H135").Value = 803.9167
$ws.Range("I135").Value = 636.2
$ws.Range("K135").Value = 5725.8
$ws.Range("M135").Value = -3190.8
$ws.Range("H138").Value = 3034.9092
$ws.Range("I138").Value = 972.1667
$ws.Range("J138").Value = 5510.2
$ws.Range("K138").Value = 2916.5001
$ws.Range("L138").Value = 16530.6
$ws.Range("M138").Value = 2223.4999
$ws.Range("N138").Value = -26810.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2550.0908
$ws.Range("I2").Value = 4183.3335
$ws.Range("K2").Value = 4183.3335
$ws.Range("M2").Value = -4070.3335
$ws.Range("H116").Value = 2550.0908
$ws.Range("I116").Value = 4183.3335
$ws.Range("K116").Value = 4183.3335
$ws.Range("M116").Value = -1889.3335
$ws.Range("H139").Value = 93499.5
$ws.Range("J139").Value = 93499.5
$ws.Range("L139").Value = 93499.5
$ws.Range("N139").Value = -103779.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2550.0908
$ws.Range("I3").Value = 4183.3335
$ws.Range("K3").Value = 4183.3335
$ws.Range("M3").Value = -4069.3335
$ws.Range("H81").Value = 33744.555
$ws.Range("J81").Value = 33744.555
$ws.Range("L81").Value = 33744.555
$ws.Range("N81").Value = -35866.555
$ws.Range("H84").Value = 33744.555
$ws.Range("J84").Value = 33744.555
$ws.Range("L84").Value = 101233.665
$ws.Range("N84").Value = -111841.665
$ws.Range("H135").Value = 177500
$ws.Range("J135").Value = 177500
$ws.Range("L135").Value = 177500
$ws.Range("N135").Value = -187640

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1735.6
$ws.Range("I22").Value = 769.8333
$ws.Range("J22").Value = 3184.25
$ws.Range("K22").Value = 769.8333
$ws.Range("L22").Value = 3184.25
$ws.Range("M22").Value = -419.8333
$ws.Range("N22").Value = -3884.25
$ws.Range("H31").Value = 6031.7617
$ws.Range("I31").Value = 1944.5
$ws.Range("K31").Value = 1944.5
$ws.Range("M31").Value = -1649.5
$ws.Range("H34").Value = 6031.7617
$ws.Range("I34").Value = 1944.5
$ws.Range("K34").Value = 1944.5
$ws.Range("M34").Value = -1742.5
$ws.Range("H62").Value = 1250
$ws.Range("I62").Value = 1250
$ws.Range("K62").Value = 1250
$ws.Range("M62").Value = -626
$ws.Range("H65").Value = 1250
$ws.Range("I65").Value = 1250
$ws.Range("K65").Value = 6250
$ws.Range("M65").Value = -3130
$ws.Range("H86").Value = 4666.3335
$ws.Range("J86").Value = 4499
$ws.Range("L86").Value = 4499
$ws.Range("N86").Value = -6745
$ws.Range("H89").Value = 4666.3335
$ws.Range("J89").Value = 4499
$ws.Range("L89").Value = 22495
$ws.Range("N89").Value = -33727
$ws.Range("H99").Value = 2469.6667
$ws.Range("I99").Value = 2463.8
$ws.Range("J99").Value = 2477
$ws.Range("K99").Value = 2463.8
$ws.Range("L99").Value = 2477
$ws.Range("M99").Value = -965.8000000000002
$ws.Range("N99").Value = -5473
$ws.Range("H126").Value = 2469.6667
$ws.Range("I126").Value = 2463.8
$ws.Range("J126").Value = 2477
$ws.Range("K126").Value = 7391.400000000001
$ws.Range("L126").Value = 7431
$ws.Range("M126").Value = -4921.400000000001
$ws.Range("N126").Value = -12371
$ws.Range("H134").Value = 1859.2222
$ws.Range("J134").Value = 1166.3334
$ws.Range("L134").Value = 3499.0002
$ws.Range("N134").Value = -8569.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 421020.2
$ws.Range("I4").Value = 501275.25
$ws.Range("K4").Value = 1503825.75
$ws.Range("M4").Value = -1503713.75
$ws.Range("H55").Value = 3938.9092
$ws.Range("J55").Value = 5304.125
$ws.Range("L55").Value = 15912.375
$ws.Range("N55").Value = -16266.375
$ws.Range("H60").Value = 1307.8125
$ws.Range("H68").Value = 392
$ws.Range("I68").Value = 392
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1176
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -365
$ws.Range("H71").Value = 392
$ws.Range("I71").Value = 392
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 3528
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 528
$ws.Range("H140").Value = 2247.6667
$ws.Range("I140").Value = 1638.25
$ws.Range("J140").Value = 4685.3335
$ws.Range("K140").Value = 4914.75
$ws.Range("L140").Value = 14056.0005
$ws.Range("M140").Value = 265.25
$ws.Range("N140").Value = -24416.0005
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 54289.5
$ws.Range("J57").Value = 54289.5
$ws.Range("L57").Value = 54289.5
$ws.Range("N57").Value = -55929.5
$ws.Range("H80").Value = 4164.75
$ws.Range("J80").Value = 4084.5
$ws.Range("L80").Value = 4084.5
$ws.Range("N80").Value = -6080.5
$ws.Range("H83").Value = 4164.75
$ws.Range("J83").Value = 4084.5
$ws.Range("L83").Value = 20422.5
$ws.Range("N83").Value = -30406.5
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("H107").Value = 650
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 300
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1620
$ws.Range("N107").Value = -4840
$ws.Range("H122").Value = 2719.2222
$ws.Range("I122").Value = 1439.4
$ws.Range("J122").Value = 4319
$ws.Range("K122").Value = 4318.200000000001
$ws.Range("L122").Value = 12957
$ws.Range("M122").Value = -1868.200000000001
$ws.Range("N122").Value = -17857
$ws.Range("N100").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1335.6666
$ws.Range("I7").Value = 1335.6666
$ws.Range("K7").Value = 1335.6666
$ws.Range("M7").Value = -1223.6666
$ws.Range("H22").Value = 1061.25
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1089.091
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1089.091
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1679.091
$ws.Range("H27").Value = 1061.25
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1089.091
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1089.091
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1303.091
$ws.Range("H46").Value = 5991.5
$ws.Range("I46").Value = 4753.6924
$ws.Range("J46").Value = 7454.364
$ws.Range("K46").Value = 4753.6924
$ws.Range("L46").Value = 7454.364
$ws.Range("M46").Value = -4565.6924
$ws.Range("N46").Value = -7830.364
$ws.Range("H126").Value = 1335.6666
$ws.Range("I126").Value = 1335.6666
$ws.Range("K126").Value = 4006.9998
$ws.Range("M126").Value = -1536.9998
$ws.Range("H132").Value = 3249.5
$ws.Range("I132").Value = 3249.5
$ws.Range("K132").Value = 9748.5
$ws.Range("M132").Value = -7218.5
$ws.Range("H136").Value = 2665.8667
$ws.Range("I136").Value = 2537.6155
$ws.Range("K136").Value = 7612.8465
$ws.Range("M136").Value = -5062.8465
